$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.752.24"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "

$ws.Cells.Item(3, 4).Value = "2.678.14"
$ws.Cells.Item(3, 5).Value = "  -0.64%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "600.42"
$ws.Cells.Item(5, 5).Value = "  -1.19%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "156.90"
$ws.Cells.Item(6, 5).Value = "  -0.79%  "

$ws.Cells.Item(8, 5).Value = "  +5.87%  "

$ws.Cells.Item(9, 5).Value = "  +5.72%  "

$ws.Cells.Item(10, 5).Value = "  -0.22%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.89"
$ws.Cells.Item(11, 5).Value = "  -2.73%  "

$ws.Cells.Item(12, 5).Value = "  -0.02%  "

$ws.Cells.Item(13, 5).Value = "  -2.71%  "

$ws.Cells.Item(14, 5).Value = "  -2.32%  "

$ws.Cells.Item(15, 4).Value = "3.158.36"
$ws.Cells.Item(15, 5).Value = "  -0.76%  "

$ws.Cells.Item(16, 4).Value = "66.320.84"
$ws.Cells.Item(16, 5).Value = "  +0.94%  "

$ws.Cells.Item(17, 4).Value = "2.670.24"
$ws.Cells.Item(17, 5).Value = "  -0.93%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "12.94"
$ws.Cells.Item(18, 5).Value = "  +2.15%  "

$ws.Cells.Item(19, 5).Value = "  -1.40%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.60"
$ws.Cells.Item(20, 5).Value = "  +1.02%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "352.61"
$ws.Cells.Item(21, 5).Value = "  -2.01%  "

$ws.Cells.Item(22, 5).Value = "  +0.05%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "69.88"
$ws.Cells.Item(23, 5).Value = "  -1.11%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.0000112"
$ws.Cells.Item(24, 5).Value = "  +4.74%  "

$ws.Cells.Item(25, 5).Value = "  -1.97%  "

$ws.Cells.Item(26, 5).Value = "  +0.25%  "

$ws.Cells.Item(27, 5).Value = "  -2.95%  "

$ws.Cells.Item(28, 5).Value = "  -5.06%  "

$ws.Cells.Item(29, 5).Value = "  -4.01%  "

$ws.Cells.Item(30, 5).Value = "  -0.05%  "

$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.14"
$ws.Cells.Item(31, 5).Value = "  -2.41%  "

$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "529.88"
$ws.Cells.Item(32, 5).Value = "  -2.98%  "

$ws.Cells.Item(33, 5).Value = "  -0.83%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.47"
$ws.Cells.Item(34, 5).Value = "  -3.69%  "

$ws.Cells.Item(35, 5).Value = "  +2.11%  "

$ws.Cells.Item(36, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.424"
$ws.Cells.Item(36, 5).Value = "  -1.90%  "

$ws.Cells.Item(37, 2).Value = "EthereumClassic"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "20.67"
$ws.Cells.Item(37, 5).Value = "  -0.65%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.999"
$ws.Cells.Item(38, 5).Value = "  +0.05%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "157.96"
$ws.Cells.Item(39, 5).Value = "  -3.32%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.94"
$ws.Cells.Item(40, 5).Value = "  -2.27%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "164.31"
$ws.Cells.Item(42, 5).Value = "  -3.03%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "4.14"
$ws.Cells.Item(43, 5).Value = "  -1.33%  "

$ws.Cells.Item(44, 5).Value = "  +1.32%  "

$ws.Cells.Item(45, 5).Value = "  -0.58%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "22.94"
$ws.Cells.Item(46, 5).Value = "  -2.78%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.642"
$ws.Cells.Item(47, 5).Value = "  -2.59%  "

$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0258"
$ws.Cells.Item(48, 5).Value = "  -3.00%  "

$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.0₆0263"
$ws.Cells.Item(49, 5).Value = "  +15.69%  "

$ws.Cells.Item(50, 5).Value = "  +1.30%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "20.20"
